# fix: mapping name on excel table
# Rename the column headers on the first table of the sheet so that
# the "Date" / "Customer" / "Total Amount" labels map to the correct columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "UserId"
$ws.Range("C1").Value = "Total Amount"
$ws.Range("D1").Value = "Date"
